$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.829.31"
$ws.Range("E2").Value = "  +2.97%  "
$ws.Range("D3").Value = "3.474.69"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.94"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.61"
$ws.Range("E6").Value = "  +5.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "3.472.31"
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  +11.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.38"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.126"
$ws.Range("E11").Value = "  +5.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").Value = "  +2.74%  "
$ws.Range("D13").Value = "4.054.54"
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("E14").Value = "  -2.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000196"
$ws.Range("E15").Value = "  +8.75%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.46"
$ws.Range("E16").Value = "  +5.76%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "64.760.25"
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.435.03"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.47"
$ws.Range("E20").Value = "  +4.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "393.77"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.33"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.96"
$ws.Range("E23").Value = "  +5.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.547"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +26.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.62"
$ws.Range("E27").Value = "  +3.59%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.21"
$ws.Range("E30").Value = "  +11.37%  "
$ws.Range("E31").Value = "  +10.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.66"
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.78"
$ws.Range("E34").Value = "  +3.31%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.06"
$ws.Range("E36").Value = "  +5.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.50"
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.27"
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0777"
$ws.Range("E40").Value = "  +4.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.44"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").Value = "2.936.78"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0321"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.86"
$ws.Range("E44").Value = "  +5.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.45"
$ws.Range("E45").Value = "  +4.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.773"
$ws.Range("E46").Value = "  +3.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.80"
$ws.Range("E47").Value = "  +8.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.10"
$ws.Range("E48").Value = "  +5.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.28"
$ws.Range("E49").Value = "  +24.03%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.108"
$ws.Range("E50").Value = "  +5.51%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.870"
$ws.Range("E51").Value = "  +7.79%  "
